# Dummy Commit to Pull Files
# Adds a new row (row 4) of test data to Sheet1:
#   A4 = TC_003 (with a left/right thin border, like A3/A2/A1 style)
#   B4 = Ajish V K
#   C4 = Ajish
# and updates the selection to follow the new last cell (C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value2 = "TC_003"
$ws.Range("B4").Value2 = "Ajish V K"
$ws.Range("C4").Value2 = "Ajish"

# Give A4 a thin left/right border (matches the border used elsewhere in the sheet).
$ws.Range("A4").Borders.Item(7).LineStyle = 1
$ws.Range("A4").Borders.Item(10).LineStyle = 1

# Move the selection to the new bottom-right data cell.
$ws.Range("C4").Select()
